$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Client Id (A2)
$ws.Range("A2").Value = "testsand316"

# Update Candidate ID (B2)
$ws.Range("B2").Value = 23071457

# Update User Name (C2)
$ws.Range("C2").Value = "igs52"

# Update Exam Password (D2)
$ws.Range("D2").Value = 'Ek6!J$z4'

# Update First Name (F2)
$ws.Range("F2").Value = "IGS"

# Update Last Name (G2)
$ws.Range("G2").Value = "India"
